$d = $word.ActiveDocument

# 1. Update the date
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line into street + city/state/zip.
#    Shrink the existing run's text to just the street address (this keeps
#    all of its original run formatting intact), then insert a brand new
#    paragraph after it - which inherits the same paragraph/run formatting
#    - and fill it in with the city/state/zip line.
$addr = $d.Content
$addr.Find.Execute("2586 Greenrock Road, Milpitas CA 95035")
$addr.Text = "2586 Greenrock Road"
$addr.InsertParagraphAfter()

$cityLine = $d.Range($addr.End + 1, $addr.End + 1)
$cityLine.Text = "Milpitas, CA 95035"

# 3. Remove the blank paragraph right after "Board of Directors"
$range = $d.Content
$range.Find.Execute("Board of Directors")
$p = $range.Paragraphs(1)
$nextPara = $p.Next()
$nextPara.Range.Delete()

$word.ActiveDocument.Save()
